# Session_spreadsheet.xlsx edit:
# 1) Correct mis-entered Animal_name values in column A (rows 215-244) so they
#    match the animal prefix of the Session_name in column B. These were
#    accidental auto-increment typos (SC002..SC031) that should all have been
#    the correct repeating animal code.
# 2) Remove the accidentally duplicated rows 264-299 (exact duplicates of the
#    now-corrected rows 228-263).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$corrections = @{
    215 = "SC001"
    216 = "SC001"
    217 = "SC001"
    218 = "SC001"
    219 = "SC001"
    220 = "SC001"
    221 = "SC001"
    222 = "SC002"
    223 = "SC002"
    224 = "SC002"
    225 = "SC002"
    226 = "SC002"
    227 = "SC002"
    228 = "SC002"
    229 = "SC004"
    230 = "SC004"
    231 = "SC004"
    232 = "SC004"
    233 = "SC004"
    234 = "SC004"
    235 = "SC004"
    236 = "SC005"
    237 = "SC005"
    238 = "SC005"
    239 = "SC005"
    240 = "SC005"
    241 = "SC005"
    242 = "SC005"
    243 = "SC005"
    244 = "SC006"
}

foreach ($row in $corrections.Keys) {
    $ws.Range("A$row").Value = $corrections[$row]
}

# Delete the duplicated tail rows (264-299), shifting everything below up.
$ws.Rows("264:299").Delete()
